$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 150, pushing the existing rows 150-213 down to 151-214.
$ws.Rows.Item(150).Insert()

# Populate the new row 150 with a new "Berenjena" price record
# (same fixed attributes as the surrounding rows, new date/volume figures).
$ws.Cells.Item(150, 1).Value = 10
$ws.Cells.Item(150, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(150, 3).Value = "La Araucanía"
$ws.Cells.Item(150, 4).Value = 44553
$ws.Cells.Item(150, 5).Value = 9
$ws.Cells.Item(150, 6).Value = 100112001
$ws.Cells.Item(150, 7).Value = "Berenjena"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 65
$ws.Cells.Item(150, 11).Value = 10000
$ws.Cells.Item(150, 12).Value = 10000
$ws.Cells.Item(150, 13).Value = 10000
$ws.Cells.Item(150, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(150, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150, 16).Value = 167
$ws.Cells.Item(150, 17).Value = 60
$ws.Cells.Item(150, 18).Value = "Hortaliza"
